# Append two new log rows (31 and 32) to the Results sheet, mirroring the
# existing run-log row 30 (same metrics), but stamped with new timestamps.
# This models two additional M3C2 stats runs that were logged after
# enhancing outlier processing/visualization (RMS + outlier threshold
# logging, PLY export for inliers/outliers).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowTemplate = @(
    @{Col='A'; Type='str'; Value='2025-08-26 14:16:15'},
    @{Col='B'; Type='str'; Value='data\0342-0349'},
    @{Col='C'; Type='str'; Value='ref'},
    @{Col='D'; Type='num'; Value=709128},
    @{Col='E'; Type='num'; Value=0.1245588149878983},
    @{Col='F'; Type='num'; Value=0.2491176299757966},
    @{Col='G'; Type='num'; Value=95},
    @{Col='H'; Type='num'; Value=0.00013396735145136},
    @{Col='I'; Type='num'; Value=0.9998660326485487},
    @{Col='J'; Type='num'; Value=709033},
    @{Col='K'; Type='num'; Value=-151.9928180000001},
    @{Col='L'; Type='num'; Value=151.729350141048},
    @{Col='M'; Type='num'; Value=690813},
    @{Col='N'; Type='num'; Value=-1125.949142},
    @{Col='O'; Type='num'; Value=62.170601304942},
    @{Col='P'; Type='num'; Value=-0.120008},
    @{Col='Q'; Type='num'; Value=0.134913},
    @{Col='R'; Type='num'; Value=-0.0002143663524828888},
    @{Col='S'; Type='num'; Value=-0.002299},
    @{Col='T'; Type='num'; Value=0.01462855990518499},
    @{Col='U'; Type='num'; Value=0.01462698916272617},
    @{Col='V'; Type='num'; Value=0.008219008269008636},
    @{Col='W'; Type='num'; Value=0.006827373},
    @{Col='X'; Type='num'; Value=-0.043862},
    @{Col='Y'; Type='num'; Value=0.043885},
    @{Col='Z'; Type='num'; Value=-0.00162988991521584},
    @{Col='AA'; Type='num'; Value=-0.002469},
    @{Col='AB'; Type='num'; Value=0.009486637003424522},
    @{Col='AC'; Type='num'; Value=0.009345573310344375},
    @{Col='AD'; Type='num'; Value=0.006657630840762984},
    @{Col='AE'; Type='num'; Value=0.0065990526},
    @{Col='AF'; Type='num'; Value=690813},
    @{Col='AG'; Type='num'; Value=246292},
    @{Col='AH'; Type='num'; Value=444509},
    @{Col='AI'; Type='num'; Value=15929},
    @{Col='AJ'; Type='num'; Value=2291},
    @{Col='AK'; Type='num'; Value=18220},
    @{Col='AL'; Type='num'; Value=0.05345534160263447},
    @{Col='AM'; Type='num'; Value=0.04536447157847037},
    @{Col='AN'; Type='num'; Value=-0.013568},
    @{Col='AO'; Type='num'; Value=-0.006508},
    @{Col='AP'; Type='num'; Value=0.002725},
    @{Col='AQ'; Type='num'; Value=0.022282},
    @{Col='AR'; Type='num'; Value=0.009233},
    @{Col='AS'; Type='num'; Value=-0.013419},
    @{Col='AT'; Type='num'; Value=-0.006579},
    @{Col='AU'; Type='num'; Value=0.002358},
    @{Col='AV'; Type='num'; Value=0.013891},
    @{Col='AW'; Type='num'; Value=0.008937},
    @{Col='AX'; Type='num'; Value=-0.0002143663524828888},
    @{Col='AY'; Type='num'; Value=0.01462698916272617},
    @{Col='AZ'; Type='num'; Value=32673381941410.1},
    @{Col='BA'; Type='num'; Value=1.010085389383569},
    @{Col='BB'; Type='num'; Value=0.3498892426805521},
    @{Col='BC'; Type='num'; Value=-0.09252700291488647},
    @{Col='BD'; Type='num'; Value=-0.08886902402080391},
    @{Col='BE'; Type='num'; Value=1.970131977073834},
    @{Col='BF'; Type='num'; Value=11637763.60543192},
    @{Col='BG'; Type='num'; Value=2.908214737699866},
    @{Col='BH'; Type='num'; Value=16.58429873827257},
    @{Col='BI'; Type='str'; Value='data\0342-0349\python_ref_m3c2_distances.txt'},
    @{Col='BJ'; Type='str'; Value='data\0342-0349\python_ref_m3c2_params.txt'},
    @{Col='BK'; Type='num'; Value=3}
)

$newTimestamps = @('2025-08-26 14:31:24', '2025-08-26 14:36:00')
$startRow = 31

for ($i = 0; $i -lt $newTimestamps.Length; $i++) {
    $r = $startRow + $i
    foreach ($entry in $rowTemplate) {
        $colLetter = $entry.Col
        $cellRef = "$colLetter$r"
        if ($entry.Type -eq 'str') {
            if ($colLetter -eq 'A') {
                $ws.Range($cellRef).Value = $newTimestamps[$i]
            } else {
                $ws.Range($cellRef).Value = $entry.Value
            }
        } else {
            $ws.Range($cellRef).Value = $entry.Value
        }
    }
}
